$d = $word.ActiveDocument

# --- Step 1: Remove the "Meta description: ..." paragraph that follows the title ---
$metaPara = $d.Paragraphs(2)
if ($metaPara.Range.Text -like "Meta description*") {
    $metaPara.Range.Delete()
}

# --- Step 2: Insert a new bold paragraph "Play Free Fire Queen Slot - Innovative
#             Gameplay Mechanism" right before the final (italic quote) paragraph ---
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$startOfLast = $lastPara.Range.Start
$insertRange = $d.Range($startOfLast, $startOfLast)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Free Fire Queen Slot - Innovative Gameplay Mechanism</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertRange.InsertXML($xml, "Before")

# InsertXML leaves a blank spacer paragraph behind it; remove that spacer so the
# new bold paragraph sits directly before the (still) final paragraph.
$count2 = $d.Paragraphs.Count
$spacerPara = $d.Paragraphs($count2 - 1)
if ($spacerPara.Range.Text.Trim() -eq "") {
    $spacerPara.Range.Delete()
}

# --- Step 3: Replace the final paragraph's italic image-prompt text with the new
#             meta-description-style sentence, keeping its italic formatting ---
$old = '"Create a feature image for Fire Queen that showcases a happy Maya warrior with glasses in a cartoon style. The image should convey the fiery and fantastical world of the game, with elements such as horses with flaming manes and phoenixes in crystal spheres incorporated into the background. The image should be vibrant and eye-catching, highlighting the unique game mechanism of Fire Queen while also staying true to the style of online slot games."'
$new = 'Discover the exciting gameplay of Fire Queen by WMS and enjoy free play. Unique mechanism and great winning potential. Suitable for innovative slot players.'
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
